# feat(EM): add blacklist mode
# Adds two new localization rows (active blacklist / active whitelist) to
# the EN localization sheet, mirroring the existing id/JP/EN row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the id/JP columns for both new rows first, then the EN column for
# both rows, so the shared-string table is populated in the same order as
# the source edit.
$ws.Range("A116").Value = "em_ui_active_blacklist"
$ws.Range("C116").Value = "アクティブなブラックリスト"

$ws.Range("A117").Value = "em_ui_active_whitelist"
$ws.Range("C117").Value = "アクティブなホワイトリスト"

$ws.Range("D116").Value = "Active Blacklist Characters"
$ws.Range("D117").Value = "Active Whitelist Characters"

# Move the active selection to D117, matching the saved cursor position.
$ws.Range("D117").Select() | Out-Null
